$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 8.908440666666667
$ws.Range("H2").Value = 26.725322
$ws.Range("I2").Value = 0.06231272032629341
$ws.Range("J2").Value = 0.06231272032629341
$ws.Range("M2").Value = 41.31538799999999
$ws.Range("N2").Value = 123.946164
$ws.Range("O2").Value = 0.2650971681450513
$ws.Range("P2").Value = 0.2650971681450514
$ws.Range("Q2").Value = 368.0556826183119
$ws.Range("R2").Value = 3312.501143564808
$ws.Range("S2").Value = 0.01651892569791496
$ws.Range("T2").Value = 0.01651892569791496

$ws.Range("G3").Value = 8.908440666666667
$ws.Range("H3").Value = 26.725322
$ws.Range("I3").Value = 0.06231272032629341
$ws.Range("J3").Value = 0.06231272032629341
$ws.Range("O3").Value = 0.5678079613724939
$ws.Range("P3").Value = 0.567807961372494
$ws.Range("Q3").Value = 788.3333808557192
$ws.Range("R3").Value = 7095.000427701472
$ws.Range("S3").Value = 0.03538165869604702
$ws.Range("T3").Value = 0.03538165869604703

$ws.Range("G4").Value = 8.908440666666667
$ws.Range("H4").Value = 26.725322
$ws.Range("I4").Value = 0.06231272032629341
$ws.Range("J4").Value = 0.06231272032629341
$ws.Range("M4").Value = 26.041732
$ws.Range("N4").Value = 78.12519599999999
$ws.Range("O4").Value = 0.1670948704824547
$ws.Range("P4").Value = 0.1670948704824547
$ws.Range("Q4").Value = 231.9912243792346
$ws.Range("R4").Value = 2087.921019413112
$ws.Range("S4").Value = 0.01041213593233142
$ws.Range("T4").Value = 0.01041213593233142

$ws.Range("I5").Value = 0.3097346304939027
$ws.Range("J5").Value = 0.3097346304939027
$ws.Range("M5").Value = 41.31538799999999
$ws.Range("N5").Value = 123.946164
$ws.Range("O5").Value = 0.2650971681450513
$ws.Range("P5").Value = 0.2650971681450514
$ws.Range("Q5").Value = 1829.475430698872
$ws.Range("R5").Value = 16465.27887628985
$ws.Range("S5").Value = 0.08210977342038746
$ws.Range("T5").Value = 0.08210977342038747

$ws.Range("I6").Value = 0.3097346304939027
$ws.Range("J6").Value = 0.3097346304939027
$ws.Range("O6").Value = 0.5678079613724939
$ws.Range("P6").Value = 0.567807961372494
$ws.Range("S6").Value = 0.1758697891072056
$ws.Range("T6").Value = 0.1758697891072056

$ws.Range("I7").Value = 0.3097346304939027
$ws.Range("J7").Value = 0.3097346304939027
$ws.Range("M7").Value = 26.041732
$ws.Range("N7").Value = 78.12519599999999
$ws.Range("O7").Value = 0.1670948704824547
$ws.Range("P7").Value = 0.1670948704824547
$ws.Range("Q7").Value = 1153.146833979741
$ws.Range("R7").Value = 10378.32150581767
$ws.Range("S7").Value = 0.05175506796630965
$ws.Range("T7").Value = 0.05175506796630965

$ws.Range("G8").Value = 89.774269
$ws.Range("H8").Value = 269.322807
$ws.Range("I8").Value = 0.627952649179804
$ws.Range("J8").Value = 0.627952649179804
$ws.Range("M8").Value = 41.31538799999999
$ws.Range("N8").Value = 123.946164
$ws.Range("O8").Value = 0.2650971681450513
$ws.Range("P8").Value = 0.2650971681450514
$ws.Range("Q8").Value = 3709.058756151371
$ws.Range("R8").Value = 33381.52880536234
$ws.Range("S8").Value = 0.1664684690267489
$ws.Range("T8").Value = 0.1664684690267489

$ws.Range("G9").Value = 89.774269
$ws.Range("H9").Value = 269.322807
$ws.Range("I9").Value = 0.627952649179804
$ws.Range("J9").Value = 0.627952649179804
$ws.Range("O9").Value = 0.5678079613724939
$ws.Range("P9").Value = 0.567807961372494
$ws.Range("Q9").Value = 7944.381698520315
$ws.Range("R9").Value = 71499.43528668283
$ws.Range("S9").Value = 0.3565565135692413
$ws.Range("T9").Value = 0.3565565135692414

$ws.Range("G10").Value = 89.774269
$ws.Range("H10").Value = 269.322807
$ws.Range("I10").Value = 0.627952649179804
$ws.Range("J10").Value = 0.627952649179804
$ws.Range("M10").Value = 26.041732
$ws.Range("N10").Value = 78.12519599999999
$ws.Range("O10").Value = 0.1670948704824547
$ws.Range("P10").Value = 0.1670948704824547
$ws.Range("Q10").Value = 2337.877453793908
$ws.Range("R10").Value = 21040.89708414517
$ws.Range("S10").Value = 0.1049276665838137
$ws.Range("T10").Value = 0.1049276665838137
